$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal text, preventing Excel from
# auto-converting numeric-looking strings (e.g. "1.000") into numbers,
# then restore the cell style so no stray formatting is introduced.
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.331.62"
Set-TextValue $ws.Range("E2") "  +0.28%  "
Set-TextValue $ws.Range("D3") "1.876.98"
Set-TextValue $ws.Range("E3") "  +0.34%  "
Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("E4") "  +0.04%  "
Set-TextValue $ws.Range("D5") "0.7119"
Set-TextValue $ws.Range("E5") "  +0.22%  "
Set-TextValue $ws.Range("D6") "242.13"
Set-TextValue $ws.Range("E6") "  +0.26%  "
Set-TextValue $ws.Range("E7") "  -0.02%  "
Set-TextValue $ws.Range("D8") "0.07873"
Set-TextValue $ws.Range("E8") "  +2.38%  "
Set-TextValue $ws.Range("D9") "0.3127"
Set-TextValue $ws.Range("E9") "  +0.75%  "
Set-TextValue $ws.Range("D10") "25.28"
Set-TextValue $ws.Range("E10") "  +1.05%  "
Set-TextValue $ws.Range("D11") "0.08388"
Set-TextValue $ws.Range("E11") "  +0.07%  "
Set-TextValue $ws.Range("D12") "1.870.81"
Set-TextValue $ws.Range("E12") "  -0.19%  "
Set-TextValue $ws.Range("D13") "5.249"
Set-TextValue $ws.Range("E13") "  +0.86%  "
Set-TextValue $ws.Range("D14") "0.7183"
Set-TextValue $ws.Range("E14") "  +1.20%  "
Set-TextValue $ws.Range("E15") "  +0.27%  "
Set-TextValue $ws.Range("D16") "6.206"
Set-TextValue $ws.Range("E16") "  +4.62%  "
Set-TextValue $ws.Range("D17") "0.000008354"
Set-TextValue $ws.Range("E17") "  +0.84%  "
Set-TextValue $ws.Range("D18") "29.330.03"
Set-TextValue $ws.Range("E18") "  +0.22%  "
Set-TextValue $ws.Range("D19") "241.08"
Set-TextValue $ws.Range("E19") "  -0.46%  "
Set-TextValue $ws.Range("D20") "13.24"
Set-TextValue $ws.Range("E20") "  +0.59%  "
Set-TextValue $ws.Range("D21") "2.121.08"
Set-TextValue $ws.Range("E21") "  -0.52%  "
Set-TextValue $ws.Range("D22") "0.9998"
Set-TextValue $ws.Range("E22") "  -0.01%  "
Set-TextValue $ws.Range("D23") "7.791"
Set-TextValue $ws.Range("E23") "  -0.46%  "
Set-TextValue $ws.Range("E24") "  +0.03%  "
Set-TextValue $ws.Range("E25") "  -2.08%  "
Set-TextValue $ws.Range("B26") "Cosmos"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D26") "9.060"
Set-TextValue $ws.Range("E26") "  +0.64%  "
Set-TextValue $ws.Range("B27") "Monero"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D27") "162.77"
Set-TextValue $ws.Range("E27") "  -0.22%  "
Set-TextValue $ws.Range("D28") "18.56"
Set-TextValue $ws.Range("E28") "  +0.51%  "
Set-TextValue $ws.Range("D29") "1.506"
Set-TextValue $ws.Range("E29") "  +0.15%  "
Set-TextValue $ws.Range("D30") "4.421"
Set-TextValue $ws.Range("E30") "  +0.38%  "
Set-TextValue $ws.Range("E31") "  +0.55%  "
Set-TextValue $ws.Range("E32") "  -4.54%  "
Set-TextValue $ws.Range("E33") "  +2.49%  "
Set-TextValue $ws.Range("D34") "1.950"
Set-TextValue $ws.Range("E34") "  +1.44%  "
Set-TextValue $ws.Range("D35") "1.180"
Set-TextValue $ws.Range("E35") "  +0.84%  "
Set-TextValue $ws.Range("D36") "0.7479"
Set-TextValue $ws.Range("E36") "  -0.39%  "
Set-TextValue $ws.Range("D37") "2.689"
Set-TextValue $ws.Range("E37") "  +0.28%  "
Set-TextValue $ws.Range("D38") "1.306.04"
Set-TextValue $ws.Range("E38") "  +12.83%  "
Set-TextValue $ws.Range("D39") "0.01882"
Set-TextValue $ws.Range("E39") "  +1.04%  "
Set-TextValue $ws.Range("D40") "2.738"
Set-TextValue $ws.Range("E40") "  +0.88%  "
Set-TextValue $ws.Range("D41") "6.519"
Set-TextValue $ws.Range("E41") "  +2.56%  "
Set-TextValue $ws.Range("B42") "TrustWalletToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "0.8953"
Set-TextValue $ws.Range("E42") "  +1.02%  "
Set-TextValue $ws.Range("B43") "Quant"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D43") "110.67"
Set-TextValue $ws.Range("E43") "  +6.38%  "
Set-TextValue $ws.Range("D44") "73.15"
Set-TextValue $ws.Range("E44") "  +0.14%  "
Set-TextValue $ws.Range("E45") "  +10.17%  "
Set-TextValue $ws.Range("D46") "0.9997"
Set-TextValue $ws.Range("E46") "  +0.00%  "
Set-TextValue $ws.Range("D47") "2.027.79"
Set-TextValue $ws.Range("E47") "  +0.11%  "
Set-TextValue $ws.Range("D48") "1.803"
Set-TextValue $ws.Range("E48") "  +0.52%  "
Set-TextValue $ws.Range("D49") "0.5191"
Set-TextValue $ws.Range("E49") "  +0.01%  "
Set-TextValue $ws.Range("D50") "9.464"
Set-TextValue $ws.Range("E50") "  +1.03%  "
Set-TextValue $ws.Range("D51") "0.4359"
Set-TextValue $ws.Range("E51") "  +1.60%  "
